$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tests")
$ws2 = $wb.Worksheets.Item("DataProviderTests")

# --- Sheet "Tests": bump D6 from 0 to 2, move the selection there ---
$ws1.Range("D6").Value = 2
[void]$ws1.Range("D6").Select()

# --- Sheet "DataProviderTests": append 5 new data-provider rows (13-17) ---
$newRows = @(
    @("navigateToTabletsPageViaHamburgerMenuTest", "yes", "chrome"),
    @("navigateToTabletsPageViaHamburgerMenuTest", "yes", "firefox"),
    @("navigateToTabletsPageViaHamburgerMenuTest", "yes", "MicrosoftEdge"),
    @("navigateToTabletsPageViaHamburgerMenuTest", "yes", "chrome"),
    @("navigateToTabletsPageViaHamburgerMenuTest", "yes", "firefox")
)

$row = 13
foreach ($data in $newRows) {
    $ws2.Cells.Item($row, 1).Value = $data[0]
    $ws2.Cells.Item($row, 2).Value = $data[1]
    $ws2.Cells.Item($row, 3).Value = $data[2]
    $row++
}

[void]$ws2.Range("B10:B17").Select()

# --- Switch the active tab back to "Tests" ---
[void]$ws1.Activate()
